$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44616
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16500
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("S2").Value = 917
$ws.Range("T2").Value = 18
$ws.Range("D3").Value = 44263
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 250
$ws.Range("Q3").Value = '$/caja 18 kilos'
$ws.Range("D4").Value = 44489
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 26000
$ws.Range("O4").Value = 27000
$ws.Range("P4").Value = 26500
$ws.Range("S4").Value = 1472
$ws.Range("D5").Value = 44307
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("S5").Value = 1083
$ws.Range("D6").Value = 44323
$ws.Range("M6").Value = 270
$ws.Range("N6").Value = 21000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 21500
$ws.Range("S6").Value = 1194
$ws.Range("D8").Value = 44487
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("S8").Value = 1450
$ws.Range("T8").Value = 10
$ws.Range("D9").Value = 44784
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 300
$ws.Range("D11").Value = 45002
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 300
$ws.Range("D12").Value = 44614
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20500
$ws.Range("S12").Value = 1139
$ws.Range("D13").Value = 44991
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 24500
$ws.Range("Q13").Value = '$/bandeja 18 kilos'
$ws.Range("S13").Value = 1361
$ws.Range("T13").Value = 18
$ws.Range("D14").Value = 44819
$ws.Range("M14").Value = 300
$ws.Range("Q14").Value = '$/bandeja 10 kilos'
$ws.Range("S14").Value = 1750
$ws.Range("T14").Value = 10
$ws.Range("D15").Value = 44629
$ws.Range("L15").Value = 'Segunda'
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 17500
$ws.Range("S15").Value = 972
$ws.Range("D16").Value = 45034
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 26000
$ws.Range("P16").Value = 25600
$ws.Range("Q16").Value = '$/bandeja 18 kilos'
$ws.Range("S16").Value = 1422
$ws.Range("T16").Value = 18
$ws.Range("D17").Value = 44291
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 17500
$ws.Range("S17").Value = 972
$ws.Range("D18").Value = 44673
$ws.Range("L18").Value = 'Especial'
$ws.Range("M18").Value = 400
$ws.Range("D19").Value = 45043
$ws.Range("L19").Value = 'Segunda'
$ws.Range("N19").Value = 21000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 21500
$ws.Range("Q19").Value = '$/bandeja 18 kilos'
$ws.Range("S19").Value = 1194
$ws.Range("T19").Value = 18
$ws.Range("D20").Value = 44491
$ws.Range("L20").Value = 'Primera'
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 14500
$ws.Range("Q20").Value = '$/bandeja 10 kilos'
$ws.Range("S20").Value = 1450
$ws.Range("T20").Value = 10
$ws.Range("D21").Value = 44706
$ws.Range("M21").Value = 400
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 9500
$ws.Range("Q21").Value = '$/bandeja 10 kilos'
$ws.Range("S21").Value = 950
$ws.Range("T21").Value = 10
$ws.Range("D22").Value = 44656
$ws.Range("M22").Value = 270
$ws.Range("N22").Value = 19000
$ws.Range("O22").Value = 20000
$ws.Range("P22").Value = 19500
$ws.Range("Q22").Value = '$/bandeja 18 kilos'
$ws.Range("S22").Value = 1083
$ws.Range("D23").Value = 44789
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 19000
$ws.Range("O23").Value = 20000
$ws.Range("P23").Value = 19500
$ws.Range("S23").Value = 1083
